$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Anatomy"

# Update row 2 with the data previously held in row 3.
# A2 must remain text (numbers-stored-as-text), not a numeric value,
# matching the original "Student ID" column convention - leading
# apostrophe forces Excel to store it as text.
$ws.Range("A2").Value = "'111111"
$ws.Range("D2").Value = "08:27:12"

# Delete row 3 (shifts nothing below it, it's the last row)
$ws.Rows.Item(3).Delete()
